$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for the movable columns (A,B,D,E,F,G,H,P,Q,R,AW,AX) per row 2-15,
# plus whether the row should carry the Z/AB ("00:00" start/end time) columns.

$rows = @{
  2  = @{ A="112243589"; B="77636"; D="NT"; E="6425";   F="Garnlav";          G="Alectoria sarmentosa";  H="(Ach.) Ach.";              P="Väst Värsjön, Vrm"; Q="404762";             R="6707097";            AW="Dick Östberg";      AX="Dick Östberg";      ZAB=$false }
  3  = @{ A="112243563"; B="89559"; D="NT"; E="5442";   F="Tallticka";        G="Porodaedalea pini";     H="(Brot.) Murrill";          P="Väst Värsjön, Vrm"; Q="404744";             R="6707084";            AW="Dick Östberg";      AX="Dick Östberg";      ZAB=$false }
  4  = @{ A="112243600"; B="77636"; D="NT"; E="6425";   F="Garnlav";          G="Alectoria sarmentosa";  H="(Ach.) Ach.";              P="Väst Värsjön, Vrm"; Q="404725";             R="6707036";            AW="Dick Östberg";      AX="Dick Östberg";      ZAB=$false }
  5  = @{ A="112243594"; B="77636"; D="NT"; E="6425";   F="Garnlav";          G="Alectoria sarmentosa";  H="(Ach.) Ach.";              P="Väst Värsjön, Vrm"; Q="404742";             R="6706992";            AW="Dick Östberg";      AX="Dick Östberg";      ZAB=$false }
  6  = @{ A="112243569"; B="77636"; D="NT"; E="6425";   F="Garnlav";          G="Alectoria sarmentosa";  H="(Ach.) Ach.";              P="Väst Värsjön, Vrm"; Q="404751";             R="6707073";            AW="Dick Östberg";      AX="Dick Östberg";      ZAB=$false }
  7  = @{ A="111525235"; B="77515"; D="NT"; E="6425";   F="Garnlav";          G="Alectoria sarmentosa";  H="(Ach.) Ach.";              P="Värsjön, väst, Vrm"; Q="404485.2245768273";  R="6706757.647421388";  AW="Helena Malmestrand"; AX="Helena Malmestrand"; ZAB=$true }
  8  = @{ A="111525233"; B="77515"; D="NT"; E="6425";   F="Garnlav";          G="Alectoria sarmentosa";  H="(Ach.) Ach.";              P="Värsjön, väst, Vrm"; Q="404540.9329893424";  R="6706716.233959051";  AW="Helena Malmestrand"; AX="Helena Malmestrand"; ZAB=$true }
  9  = @{ A="111525238"; B="77515"; D="NT"; E="6425";   F="Garnlav";          G="Alectoria sarmentosa";  H="(Ach.) Ach.";              P="Värsjön, väst, Vrm"; Q="404495.4563026094";  R="6706677.491168984";  AW="Helena Malmestrand"; AX="Helena Malmestrand"; ZAB=$true }
  10 = @{ A="112243588"; B="77636"; D="NT"; E="6425";   F="Garnlav";          G="Alectoria sarmentosa";  H="(Ach.) Ach.";              P="Väst Värsjön, Vrm"; Q="404452";             R="6706739";            AW="Dick Östberg";      AX="Dick Östberg";      ZAB=$false }
  11 = @{ A="112243565"; B="90800"; D="LC"; E="4364";   F="Dropptaggsvamp";   G="Hydnellum ferrugineum"; H="(Fr.:Fr.) P. Karst.";      P="Väst Värsjön, Vrm"; Q="404459";             R="6706753";            AW="Dick Östberg";      AX="Dick Östberg";      ZAB=$false }
  12 = @{ A="112243573"; B="77636"; D="NT"; E="6425";   F="Garnlav";          G="Alectoria sarmentosa";  H="(Ach.) Ach.";              P="Väst Värsjön, Vrm"; Q="404477";             R="6706766";            AW="Dick Östberg";      AX="Dick Östberg";      ZAB=$false }
  13 = @{ A="111525223"; B="78107"; D="NT"; E="6453";   F="Vedskivlav";       G="Hertelidea botryosa";   H="(Fr.) Printzen & Kantvilas"; P="Värsjön, väst, Vrm"; Q="404637.0659126193"; R="6706784.214121711"; AW="Helena Malmestrand"; AX="Helena Malmestrand"; ZAB=$true }
  14 = @{ A="111525224"; B="77268"; D="NT"; E="228912"; F="Mörk kolflarnlav"; G="Carbonicola myrmecina"; H="(Ach.) Bendiksby & Timdal"; P="Värsjön, väst, Vrm"; Q="404619.9854206198"; R="6706773.322858612"; AW="Helena Malmestrand"; AX="Helena Malmestrand"; ZAB=$true }
  15 = @{ A="111525226"; B="77515"; D="NT"; E="6425";   F="Garnlav";          G="Alectoria sarmentosa";  H="(Ach.) Ach.";              P="Värsjön, väst, Vrm"; Q="404616.9589749529";  R="6706770.937089294";  AW="Helena Malmestrand"; AX="Helena Malmestrand"; ZAB=$true }
}

foreach ($r in $rows.Keys) {
  $vals = $rows[$r]

  $ws.Range("A$r").Value = [double]$vals.A
  $ws.Range("B$r").Value = [double]$vals.B
  $ws.Range("D$r").Value = $vals.D
  $ws.Range("E$r").Value = [double]$vals.E
  $ws.Range("F$r").Value = $vals.F
  $ws.Range("G$r").Value = $vals.G
  $ws.Range("H$r").Value = $vals.H
  $ws.Range("P$r").Value = $vals.P
  $ws.Range("Q$r").Value = [double]$vals.Q
  $ws.Range("R$r").Value = [double]$vals.R
  $ws.Range("AW$r").Value = $vals.AW
  $ws.Range("AX$r").Value = $vals.AX

  if ($vals.ZAB) {
    $ws.Range("Z$r").Value = "00:00"
    $ws.Range("AB$r").Value = "00:00"
  } else {
    $ws.Range("Z$r").Value = ""
    $ws.Range("AB$r").Value = ""
  }
}
